# ============================================================================
# template_vm.xlsx edit: "update user log, add replikasi vm, dll"
# Rebuilds Sheet1 header/data layout: drops Host/Processor columns, adds
# Core per Socket, Jumlah Socket, Aplikasi ID, Masa Aktif, Memo VM,
# Status Replikasi and a formatted Note column; adds dropdown validations
# and a date value; refreshes workbook session metadata.
# ============================================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet so the old Host/Processor columns and stale
# styles/shared-strings don't linger.
$ws.Cells.Clear()

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Rows(1).RowHeight = 135

$headers = @{
    "A1" = "Cluster ID";
    "B1" = "OS ID";
    "C1" = "Nama VM";
    "D1" = "IP Address";
    "E1" = "Hostname";
    "F1" = "Disk (Gb)";
    "G1" = "Memory (Gb)";
    "H1" = "Core per Socket";
    "I1" = "Jumlah Socket";
    "J1" = "Jenis Server";
    "K1" = "Lisence";
    "L1" = "Aplikasi ID";
    "M1" = "Masa Aktif";
    "N1" = "Memo VM";
    "O1" = "Status Replikasi";
}
foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# Blue/bold header cells (Cluster ID, OS ID, Aplikasi ID)
foreach ($addr in @("A1","B1","L1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.ThemeColor = 5
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Black/bold header cells
foreach ($addr in @("C1","D1","E1","F1","G1","H1","I1","K1","M1","N1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Red/bold header cells (Jenis Server, Status Replikasi)
foreach ($addr in @("J1","O1")) {
    $cell = $ws.Range($addr)
    $cell.Font.Bold = $true
    $cell.Font.Color = 255
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Note column header (P1) - multi-run rich text, red, left/wrap aligned.
$noteText = "*Note :`n- Kolom isian tidak boleh kosong`n- Input Cluster ID, OS ID dan Aplikasi ID sesuai dengan ID yang terdaftar`n- Jenis Server dan Status Replikasi diisi dengan cara memilih dropdown yang sudah tersedia`n- Kolom Memo VM, isi dengan kata `"kosong`" (tanpa tanda petik)`n- Pastikan seluruh kolom yang kosong dihapus, dengan cara blok kolom yang kosong kemudian klik kanan -> pilih delete  -> ok"
$ws.Range("P1").Value = $noteText

$p1 = $ws.Range("P1")
$p1.Font.Color = 255
$p1.HorizontalAlignment = -4131
$p1.VerticalAlignment = -4108
$p1.WrapText = $true

# "dihapus" run -> bold red
$boldStart = 306
$boldLen = 7
$p1.Characters($boldStart, $boldLen).Font.Bold = $true
$p1.Characters($boldStart, $boldLen).Font.Color = 255

# trailing run (", dengan cara ... -> ok") -> red, not bold
$tailStart = $boldStart + $boldLen
$tailLen = 79
$p1.Characters($tailStart, $tailLen).Font.Color = 255
$p1.Characters($tailStart, $tailLen).Font.Bold = $false

# ---------------------------------------------------------------------------
# Data row (row 2)
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "VM1"
$ws.Range("D2").Value = "172.3.1.1"
$ws.Range("E2").Value = "HOSTNAMEVM1"
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = "WEB"
$ws.Range("K2").Value = "LVM1"
$ws.Range("L2").Value = 1

$ws.Range("M2").Value = 45387
$ws.Range("M2").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("N2").Value = "kosong"
$ws.Range("N2").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("O2").Value = "site_recovery_manajemen"

# ---------------------------------------------------------------------------
# Data validations
# ---------------------------------------------------------------------------
$jv = $ws.Range("J2").Validation
$jv.Add(3, 1, 1, '"WEB,APP,DB,MNGMT,DMZ,DEV"')
$jv.IgnoreBlank = $false

$ov = $ws.Range("O2").Validation
$ov.Add(3, 1, 1, '"database_replikasi,site_recovery_manajemen,belum_replikasi"')
$ov.IgnoreBlank = $false

# ---------------------------------------------------------------------------
# Column widths (best fit, approximated)
# ---------------------------------------------------------------------------
$ws.Columns(1).ColumnWidth = 8.67
$ws.Columns(2).ColumnWidth = 4.83
$ws.Columns(3).ColumnWidth = 9.71
$ws.Columns(4).ColumnWidth = 9.5
$ws.Columns(5).ColumnWidth = 14.5
$ws.Columns(6).ColumnWidth = 8.33
$ws.Columns(7).ColumnWidth = 12.33
$ws.Columns(8).ColumnWidth = 14.17
$ws.Columns(9).ColumnWidth = 14.83
$ws.Columns(10).ColumnWidth = 10.67
$ws.Columns(11).ColumnWidth = 6.67
$ws.Columns(12).ColumnWidth = 9.5
$ws.Columns(13).ColumnWidth = 9.83
$ws.Columns(14).ColumnWidth = 9.83
$ws.Columns(15).ColumnWidth = 29.17
$ws.Columns(16).ColumnWidth = 62.83

$ws.Columns(17).ColumnWidth = -0.8333333333333334
$ws.Columns(17).Hidden = $true

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("A3:P1048576").Select()

# ---------------------------------------------------------------------------
# Workbook session metadata
# ---------------------------------------------------------------------------
$wb.AbsPath = "E:\template_xls\"
